$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.266.34'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.915.17'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7353'
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.23'
$ws.Range("E6").Value = '  -2.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.007'
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3120'
$ws.Range("E8").Value = '  -2.95%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.13'
$ws.Range("E9").Value = '  -2.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06903'
$ws.Range("E10").Value = '  -2.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07996'
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7650'
$ws.Range("E12").Value = '  -2.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.901.05'
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.276'
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.04'
$ws.Range("E15").Value = '  -3.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.226.65'
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.13'
$ws.Range("E17").Value = '  -3.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '245.34'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.831'
$ws.Range("E19").Value = '  +1.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007790'
$ws.Range("E20").Value = '  -3.09%  '
$ws.Range("E21").Value = '  +0.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.155.27'
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.008'
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.566'
$ws.Range("E24").Value = '  -3.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.333'
$ws.Range("E25").Value = '  -2.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.22'
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.79'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1272'
$ws.Range("E28").Value = '  -3.85%  '
$ws.Range("E29").Value = '  -7.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.369'
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.542'
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.319'
$ws.Range("E32").Value = '  -2.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.040'
$ws.Range("E33").Value = '  -2.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05141'
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.289'
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7398'
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.782'
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01926'
$ws.Range("E38").Value = '  -2.63%  '
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.399'
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '75.38'
$ws.Range("E41").Value = '  -3.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4430'
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.926'
$ws.Range("E43").Value = '  -2.98%  '
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8352'
$ws.Range("E45").Value = '  -0.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.79'
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.543'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.779'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.92'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '950.55'
$ws.Range("E50").Value = '  -2.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1198'
$ws.Range("E51").Value = '  +3.98%  '
